$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B data entry (soil analysis sample values) ---
$ws.Range("B2").Value = "Soil"
$ws.Range("B3").Value = "test"
$ws.Range("B4").Value = 1
$ws.Range("B5").Value = "Yes"
$ws.Range("B6").Value = 12

# B7 holds a date (Sample_Date). Copy the number-formatting/font/alignment
# from an existing centered cell (A11) so the new style shares the same
# font/alignment, then apply the date number format and finally set the
# date value (ordering matters: format before value avoids creating a
# spurious auxiliary number format entry).
$ws.Range("A11").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4122) | Out-Null
$ws.Range("B7").NumberFormat = "mm-dd-yy"
$ws.Range("B7").Value = (Get-Date -Year 2019 -Month 5 -Day 2 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B8").Value = 11
$ws.Range("B9").Value = "Silty_Clay"
$ws.Range("B10").Value = "0-30"
$ws.Range("B11").Value = 2
$ws.Range("B12").Value = 0
$ws.Range("B13").Value = 0
$ws.Range("B14").Value = 50
$ws.Range("B15").Value = 2
$ws.Range("B16").Value = 3
$ws.Range("B17").Value = 4

# --- Selection moves to B8 ---
$ws.Range("B8").Select() | Out-Null
